# Daily attendance processing - 2025-12-31 16:37:33
# Reorders the comma-separated names/emails in the "Recorded By" (column G)
# cells on the "Session Analysis Results" sheet, moving the first listed
# entry to the end (i.e. reversing the order of the comma-separated list).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ', '
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(', ', $reversed)
        }
    }
}
